# Adding the User Transaction feature SIQ and SRS
# - Rename "Sheet1" -> "RTM"
# - Add a new "revision record" sheet after it, documenting the V1.0 / V2.0 history
# - Make the new sheet the active tab (tabSelected)

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet -------------------------------------------------
$rtm = $wb.Worksheets.Item(1)
$rtm.Name = "RTM"

# --- Add the new "revision record" worksheet right after RTM -------------------
$rev = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rtm)
$rev.Name = "revision record"

# --- Header row ------------------------------------------------------------
$rev.Range("A1").Value = "version"
$rev.Range("B1").Value = "Date"
$rev.Range("C1").Value = "owner"
$rev.Range("D1").Value = "Description"

$header = $rev.Range("A1:D1")
$header.Font.Bold = $true
$header.Font.Size = 12
$header.Interior.ThemeColor = 8

# --- Data rows ---------------------------------------------------------------
$rev.Range("A2").Value = "V1.0"
$rev.Range("B2").Value = (Get-Date -Year 2024 -Month 8 -Day 4).Date
$rev.Range("B2").NumberFormat = "m/d/yyyy"
$rev.Range("C2").Value = "Mayar"
$rev.Range("D2").Value = "RTM  without User transaction feature"

$rev.Range("A3").Value = "V2.0"
$rev.Range("B3").Value = "19/4/2024"
$rev.Range("C3").Value = "Mayar"
$rev.Range("D3").Value = "Update RTM to contain User transaction feature"

# --- Column widths to roughly match the authored layout ------------------------
$rev.Columns.Item(1).ColumnWidth = 10.36328125
$rev.Columns.Item(2).ColumnWidth = 17
$rev.Columns.Item(3).ColumnWidth = 17.453125
$rev.Columns.Item(4).ColumnWidth = 44

# --- Make "revision record" the active/visible tab -----------------------------
$rev.Select()
$rev.Range("B4").Select()
